$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = 658607
$ws.Cells.Item(2, 5).Value = -1828
$ws.Cells.Item(2, 6).Value = -2313
$ws.Cells.Item(2, 7).Value = -4367
$ws.Cells.Item(2, 8).Value = -5372
$ws.Cells.Item(2, 9).Value = -5888
$ws.Cells.Item(2, 10).Value = 517
$ws.Cells.Item(2, 11).Value = 351013
$ws.Cells.Item(2, 12).Value = 190403
$ws.Cells.Item(2, 13).Value = 160611
$ws.Cells.Item(2, 14).Value = 150003
$ws.Cells.Item(2, 15).Value = 10607
$ws.Cells.Item(2, 16).Value = 4686
$ws.Cells.Item(2, 17).Value = 9271
$ws.Cells.Item(2, 18).Value = -25694
$ws.Cells.Item(2, 19).Value = 17099
$ws.Cells.Item(2, 20).Value = 16351
$ws.Cells.Item(2, 21).Value = -7080
$ws.Cells.Item(2, 22).Value = 109752
$ws.Cells.Item(2, 23).Value = -0.28
$ws.Cells.Item(2, 24).Value = -0.82
$ws.Cells.Item(2, 25).Value = -3.82
$ws.Cells.Item(2, 26).Value = -1.53
$ws.Cells.Item(2, 27).Value = 118.55
$ws.Cells.Item(2, 28).Value = 3100.74
$ws.Cells.Item(2, 29).Value = -6283
$ws.Cells.Item(2, 30).Value = -13.54
$ws.Cells.Item(2, 31).Value = 160961
$ws.Cells.Item(2, 32).Value = 0.53
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 92465564
$ws.Cells.Item(3, 4).Value = 483563
$ws.Cells.Item(3, 5).Value = 19796
$ws.Cells.Item(3, 6).Value = 19796
$ws.Cells.Item(3, 7).Value = 16541
$ws.Cells.Item(3, 8).Value = 8677
$ws.Cells.Item(3, 9).Value = 8149
$ws.Cells.Item(3, 10).Value = 527
$ws.Cells.Item(3, 11).Value = 313598
$ws.Cells.Item(3, 12).Value = 143135
$ws.Cells.Item(3, 13).Value = 170463
$ws.Cells.Item(3, 14).Value = 158897
$ws.Cells.Item(3, 15).Value = 11566
$ws.Cells.Item(3, 16).Value = 4686
$ws.Cells.Item(3, 17).Value = 40857
$ws.Cells.Item(3, 18).Value = -13454
$ws.Cells.Item(3, 19).Value = -27683
$ws.Cells.Item(3, 20).Value = 4488
$ws.Cells.Item(3, 21).Value = 36369
$ws.Cells.Item(3, 22).Value = 81748
$ws.Cells.Item(3, 23).Value = 4.09
$ws.Cells.Item(3, 24).Value = 1.79
$ws.Cells.Item(3, 25).Value = 5.28
$ws.Cells.Item(3, 26).Value = 2.61
$ws.Cells.Item(3, 27).Value = 83.97
$ws.Cells.Item(3, 28).Value = 3272.47
$ws.Cells.Item(3, 29).Value = 8696
$ws.Cells.Item(3, 30).Value = 14.95
$ws.Cells.Item(3, 31).Value = 170504
$ws.Cells.Item(3, 32).Value = 0.76
$ws.Cells.Item(3, 33).Value = 4800
$ws.Cells.Item(3, 34).Value = 3.69
$ws.Cells.Item(3, 35).Value = 54.9
$ws.Cells.Item(3, 36).Value = 92465564
$ws.Cells.Item(4, 4).Value = 395205
$ws.Cells.Item(4, 5).Value = 32283
$ws.Cells.Item(4, 6).Value = 32283
$ws.Cells.Item(4, 7).Value = 24234
$ws.Cells.Item(4, 8).Value = 17214
$ws.Cells.Item(4, 9).Value = 16713
$ws.Cells.Item(4, 10).Value = 501
$ws.Cells.Item(4, 11).Value = 325813
$ws.Cells.Item(4, 12).Value = 142763
$ws.Cells.Item(4, 13).Value = 183050
$ws.Cells.Item(4, 14).Value = 171032
$ws.Cells.Item(4, 15).Value = 12019
$ws.Cells.Item(4, 16).Value = 4686
$ws.Cells.Item(4, 17).Value = 36778
$ws.Cells.Item(4, 18).Value = -20561
$ws.Cells.Item(4, 19).Value = -20186
$ws.Cells.Item(4, 20).Value = 6070
$ws.Cells.Item(4, 21).Value = 30708
$ws.Cells.Item(4, 22).Value = 65697
$ws.Cells.Item(4, 23).Value = 8.17
$ws.Cells.Item(4, 24).Value = 4.36
$ws.Cells.Item(4, 25).Value = 10.13
$ws.Cells.Item(4, 26).Value = 5.38
$ws.Cells.Item(4, 27).Value = 77.98999999999999
$ws.Cells.Item(4, 28).Value = 3533.1
$ws.Cells.Item(4, 29).Value = 17834
$ws.Cells.Item(4, 30).Value = 8.210000000000001
$ws.Cells.Item(4, 31).Value = 183525
$ws.Cells.Item(4, 32).Value = 0.8
$ws.Cells.Item(4, 33).Value = 6400
$ws.Cells.Item(4, 34).Value = 4.37
$ws.Cells.Item(4, 35).Value = 35.69
$ws.Cells.Item(4, 36).Value = 92465564
$ws.Cells.Item(5, 4).Value = 461627
$ws.Cells.Item(5, 5).Value = 32218
$ws.Cells.Item(5, 6).Value = 32218
$ws.Cells.Item(5, 7).Value = 32237
$ws.Cells.Item(5, 8).Value = 21451
$ws.Cells.Item(5, 9).Value = 21038
$ws.Cells.Item(5, 10).Value = 413
$ws.Cells.Item(5, 11).Value = 342501
$ws.Cells.Item(5, 12).Value = 149408
$ws.Cells.Item(5, 13).Value = 193093
$ws.Cells.Item(5, 14).Value = 180858
$ws.Cells.Item(5, 15).Value = 12236
$ws.Cells.Item(5, 16).Value = 4686
$ws.Cells.Item(5, 17).Value = 21802
$ws.Cells.Item(5, 18).Value = -10661
$ws.Cells.Item(5, 19).Value = -16707
$ws.Cells.Item(5, 20).Value = 9384
$ws.Cells.Item(5, 21).Value = 12418
$ws.Cells.Item(5, 22).Value = 55779
$ws.Cells.Item(5, 23).Value = 6.98
$ws.Cells.Item(5, 24).Value = 4.65
$ws.Cells.Item(5, 25).Value = 11.96
$ws.Cells.Item(5, 26).Value = 6.42
$ws.Cells.Item(5, 27).Value = 77.38
$ws.Cells.Item(5, 28).Value = 3800.39
$ws.Cells.Item(5, 29).Value = 22449
$ws.Cells.Item(5, 30).Value = 9.109999999999999
$ws.Cells.Item(5, 31).Value = 194069
$ws.Cells.Item(5, 32).Value = 1.05
$ws.Cells.Item(5, 33).Value = 8000
$ws.Cells.Item(5, 34).Value = 3.91
$ws.Cells.Item(5, 35).Value = 35.44
$ws.Cells.Item(5, 36).Value = 92465564
$ws.Cells.Item(6, 4).Value = 545109
$ws.Cells.Item(6, 5).Value = 21176
$ws.Cells.Item(6, 6).Value = 21176
$ws.Cells.Item(6, 7).Value = 24024
$ws.Cells.Item(6, 8).Value = 17100
$ws.Cells.Item(6, 9).Value = 16515
$ws.Cells.Item(6, 11).Value = 360854
$ws.Cells.Item(6, 12).Value = 167574
$ws.Cells.Item(6, 13).Value = 193280
$ws.Cells.Item(6, 14).Value = 181244
$ws.Cells.Item(6, 16).Value = 4686
$ws.Cells.Item(6, 17).Value = 17281
$ws.Cells.Item(6, 18).Value = -24768
$ws.Cells.Item(6, 19).Value = 5861
$ws.Cells.Item(6, 20).Value = 12789
$ws.Cells.Item(6, 21).Value = 4491
$ws.Cells.Item(6, 22).Value = 80233
$ws.Cells.Item(6, 23).Value = 3.89
$ws.Cells.Item(6, 24).Value = 3.14
$ws.Cells.Item(6, 25).Value = 9.119999999999999
$ws.Cells.Item(6, 26).Value = 4.86
$ws.Cells.Item(6, 27).Value = 86.7
$ws.Cells.Item(6, 28).Value = 3990.04
$ws.Cells.Item(6, 29).Value = 17622
$ws.Cells.Item(6, 30).Value = 10.19
$ws.Cells.Item(6, 31).Value = 205466
$ws.Cells.Item(6, 32).Value = 0.87
$ws.Cells.Item(6, 33).Value = 8000
$ws.Cells.Item(6, 34).Value = 4.46
$ws.Cells.Item(6, 35).Value = 42.89
$ws.Cells.Item(6, 36).Value = 92465564
$ws.Cells.Item(7, 4).Value = 505947
$ws.Cells.Item(7, 5).Value = 13635
$ws.Cells.Item(7, 7).Value = 10163
$ws.Cells.Item(7, 8).Value = 7013
$ws.Cells.Item(7, 9).Value = 6405
$ws.Cells.Item(7, 11).Value = 386701
$ws.Cells.Item(7, 12).Value = 190341
$ws.Cells.Item(7, 13).Value = 196360
$ws.Cells.Item(7, 14).Value = 182083
$ws.Cells.Item(7, 16).Value = 4688
$ws.Cells.Item(7, 17).Value = 21045
$ws.Cells.Item(7, 18).Value = -31645
$ws.Cells.Item(7, 19).Value = 15093
$ws.Cells.Item(7, 20).Value = 28661
$ws.Cells.Item(7, 21).Value = -5278
$ws.Cells.Item(7, 23).Value = 2.7
$ws.Cells.Item(7, 24).Value = 1.39
$ws.Cells.Item(7, 25).Value = 3.53
$ws.Cells.Item(7, 26).Value = 1.88
$ws.Cells.Item(7, 27).Value = 96.93000000000001
$ws.Cells.Item(7, 29).Value = 6834
$ws.Cells.Item(7, 30).Value = 19.02
$ws.Cells.Item(7, 31).Value = 206416
$ws.Cells.Item(7, 32).Value = 0.63
$ws.Cells.Item(7, 33).Value = 6105
$ws.Cells.Item(7, 34).Value = 4.7
$ws.Cells.Item(7, 35).Value = 88.13
$ws.Cells.Item(8, 4).Value = 506386
$ws.Cells.Item(8, 5).Value = 17600
$ws.Cells.Item(8, 7).Value = 16579
$ws.Cells.Item(8, 8).Value = 11703
$ws.Cells.Item(8, 9).Value = 11145
$ws.Cells.Item(8, 11).Value = 402958
$ws.Cells.Item(8, 12).Value = 200116
$ws.Cells.Item(8, 13).Value = 202842
$ws.Cells.Item(8, 14).Value = 188268
$ws.Cells.Item(8, 16).Value = 4688
$ws.Cells.Item(8, 17).Value = 23037
$ws.Cells.Item(8, 18).Value = -25982
$ws.Cells.Item(8, 19).Value = 2672
$ws.Cells.Item(8, 20).Value = 25294
$ws.Cells.Item(8, 21).Value = 1398
$ws.Cells.Item(8, 23).Value = 3.48
$ws.Cells.Item(8, 24).Value = 2.31
$ws.Cells.Item(8, 25).Value = 6.02
$ws.Cells.Item(8, 26).Value = 2.96
$ws.Cells.Item(8, 27).Value = 98.66
$ws.Cells.Item(8, 29).Value = 11892
$ws.Cells.Item(8, 30).Value = 10.93
$ws.Cells.Item(8, 31).Value = 213428
$ws.Cells.Item(8, 32).Value = 0.61
$ws.Cells.Item(8, 33).Value = 6719
$ws.Cells.Item(8, 34).Value = 5.17
$ws.Cells.Item(8, 35).Value = 55.75
$ws.Cells.Item(9, 4).Value = 511958
$ws.Cells.Item(9, 5).Value = 20369
$ws.Cells.Item(9, 7).Value = 18150
$ws.Cells.Item(9, 8).Value = 12829
$ws.Cells.Item(9, 9).Value = 12213
$ws.Cells.Item(9, 11).Value = 419997
$ws.Cells.Item(9, 12).Value = 210232
$ws.Cells.Item(9, 13).Value = 209765
$ws.Cells.Item(9, 14).Value = 194917
$ws.Cells.Item(9, 16).Value = 4688
$ws.Cells.Item(9, 17).Value = 25584
$ws.Cells.Item(9, 18).Value = -27057
$ws.Cells.Item(9, 19).Value = 3372
$ws.Cells.Item(9, 20).Value = 24588
$ws.Cells.Item(9, 21).Value = 3743
$ws.Cells.Item(9, 23).Value = 3.98
$ws.Cells.Item(9, 24).Value = 2.51
$ws.Cells.Item(9, 25).Value = 6.37
$ws.Cells.Item(9, 26).Value = 3.12
$ws.Cells.Item(9, 27).Value = 100.22
$ws.Cells.Item(9, 29).Value = 13032
$ws.Cells.Item(9, 30).Value = 9.98
$ws.Cells.Item(9, 31).Value = 220965
$ws.Cells.Item(9, 32).Value = 0.59
$ws.Cells.Item(9, 33).Value = 6900
$ws.Cells.Item(9, 34).Value = 5.31
$ws.Cells.Item(9, 35).Value = 52.24
